# Generate Report for Handback
# Update the handback status timestamps / priority that changed between
# report generation runs (these cells already hold plain text values,
# so plain assignment keeps them as text rather than real dates).

$wb = $excel.ActiveWorkbook

# --- "Overview" sheet: Latest HO Xliff Generate Date (column G) ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G3").Value = "2016-08-23 22:14:08"
$wsOverview.Range("G5").Value = "2016-08-23 22:14:08"

# --- "zh-cn" sheet: Priority (E), Correspond Handoff Datetime (H),
#     Correspond Handback DateTime (K) ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Range("E3").Value = "mt"
$wsZhCn.Range("E5").Value = "mt"

$wsZhCn.Range("H3").Value = "2016-08-23 22:13:58"
$wsZhCn.Range("H5").Value = "2016-08-23 22:13:58"

$wsZhCn.Range("K3").Value = "2016-08-23 22:14:27"
$wsZhCn.Range("K5").Value = "2016-08-23 22:14:27"

# --- "de-de" sheet: Priority (E), Correspond Handoff Datetime (H),
#     Correspond Handback DateTime (K) ---
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Range("E3").Value = "mt"
$wsDeDe.Range("E5").Value = "mt"

$wsDeDe.Range("H3").Value = "2016-08-23 22:14:08"
$wsDeDe.Range("H5").Value = "2016-08-23 22:14:08"

$wsDeDe.Range("K3").Value = "2016-08-23 22:14:34"
$wsDeDe.Range("K5").Value = "2016-08-23 22:14:34"
